$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 3082.9473
$ws.Range("I111").Value = 3161.2307
$ws.Range("J111").Value = 2913.3333
$ws.Range("K111").Value = 9483.6921
$ws.Range("L111").Value = 8739.999899999999
$ws.Range("M111").Value = -6416.6921
$ws.Range("N111").Value = -14873.9999
$ws.Range("H112").Value = 1972.9166
$ws.Range("J112").Value = 2307.5
$ws.Range("L112").Value = 6922.5
$ws.Range("N112").Value = -9138.5
$ws.Range("H138").Value = 4439.519
$ws.Range("J138").Value = 4417.1816
$ws.Range("L138").Value = 13251.5448
$ws.Range("N138").Value = -23531.5448

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 74230.23
$ws.Range("I32").Value = 26484.291
$ws.Range("K32").Value = 26484.291
$ws.Range("M32").Value = -26197.291
$ws.Range("H88").Value = 1919.5
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("H91").Value = 1919.5
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("H97").Value = 43027.418
$ws.Range("I97").Value = 46619.816
$ws.Range("J97").Value = 3511
$ws.Range("K97").Value = 46619.816
$ws.Range("L97").Value = 3511
$ws.Range("M97").Value = -46123.816
$ws.Range("N97").Value = -4503
$ws.Range("H110").Value = 41755640
$ws.Range("I110").Value = 47720544
$ws.Range("J110").Value = 1321
$ws.Range("K110").Value = 47720544
$ws.Range("L110").Value = 1321
$ws.Range("M110").Value = -47718499
$ws.Range("N110").Value = -5411
$ws.Range("H132").Value = 15172060
$ws.Range("I132").Value = 17264172
$ws.Range("J132").Value = 4249
$ws.Range("K132").Value = 51792516
$ws.Range("L132").Value = 12747
$ws.Range("M132").Value = -51789986
$ws.Range("N132").Value = -17807

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 33074.234
$ws.Range("I20").Value = 44374.75
$ws.Range("J20").Value = 5953
$ws.Range("K20").Value = 44374.75
$ws.Range("L20").Value = 5953
$ws.Range("M20").Value = -44127.75
$ws.Range("N20").Value = -6447
$ws.Range("H86").Value = 87726.16
$ws.Range("I86").Value = 103390
$ws.Range("J86").Value = 1575
$ws.Range("K86").Value = 103390
$ws.Range("L86").Value = 1575
$ws.Range("M86").Value = -102267
$ws.Range("N86").Value = -3821
$ws.Range("H89").Value = 87726.16
$ws.Range("I89").Value = 103390
$ws.Range("J89").Value = 1575
$ws.Range("K89").Value = 516950
$ws.Range("L89").Value = 7875
$ws.Range("M89").Value = -511334
$ws.Range("N89").Value = -19107
$ws.Range("H94").Value = 280.13333
$ws.Range("I94").Value = 263.31708
$ws.Range("J94").Value = 452.5
$ws.Range("K94").Value = 263.31708
$ws.Range("L94").Value = 452.5
$ws.Range("M94").Value = 187.68292
$ws.Range("N94").Value = -1354.5
$ws.Range("H105").Value = 81886.32000000001
$ws.Range("I105").Value = 85301.586
$ws.Range("J105").Value = 78733.766
$ws.Range("K105").Value = 85301.586
$ws.Range("L105").Value = 78733.766
$ws.Range("M105").Value = -83554.586
$ws.Range("N105").Value = -82227.766
$ws.Range("H107").Value = 125060500
$ws.Range("I107").Value = 166743580
$ws.Range("J107").Value = 11237.5
$ws.Range("K107").Value = 166743580
$ws.Range("L107").Value = 11237.5
$ws.Range("M107").Value = -166741660
$ws.Range("N107").Value = -15077.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1160.0769
$ws.Range("I16").Value = 881.1429000000001
$ws.Range("J16").Value = 1485.5
$ws.Range("K16").Value = 881.1429000000001
$ws.Range("L16").Value = 1485.5
$ws.Range("M16").Value = -594.1429000000001
$ws.Range("N16").Value = -2059.5
$ws.Range("H58").Value = 1605.5897
$ws.Range("I58").Value = 1366.2667
$ws.Range("J58").Value = 2403.3333
$ws.Range("K58").Value = 1366.2667
$ws.Range("L58").Value = 2403.3333
$ws.Range("M58").Value = -1163.2667
$ws.Range("N58").Value = -2809.3333
$ws.Range("H74").Value = 36657
$ws.Range("J74").Value = 36657
$ws.Range("L74").Value = 36657
$ws.Range("N74").Value = -38405
$ws.Range("H77").Value = 36657
$ws.Range("J77").Value = 36657
$ws.Range("L77").Value = 109971
$ws.Range("N77").Value = -118707
$ws.Range("H86").Value = 3712.7666
$ws.Range("I86").Value = 3282.389
$ws.Range("J86").Value = 4358.3335
$ws.Range("K86").Value = 3282.389
$ws.Range("L86").Value = 4358.3335
$ws.Range("M86").Value = -2159.389
$ws.Range("N86").Value = -6604.3335
$ws.Range("H89").Value = 3712.7666
$ws.Range("I89").Value = 3282.389
$ws.Range("J89").Value = 4358.3335
$ws.Range("K89").Value = 16411.945
$ws.Range("L89").Value = 21791.6675
$ws.Range("M89").Value = -10795.945
$ws.Range("N89").Value = -33023.6675
$ws.Range("H105").Value = 998.43475
$ws.Range("I105").Value = 1102
$ws.Range("J105").Value = 705
$ws.Range("K105").Value = 1102
$ws.Range("L105").Value = 705
$ws.Range("M105").Value = 645
$ws.Range("N105").Value = -4199
$ws.Range("H113").Value = 1160.0769
$ws.Range("I113").Value = 881.1429000000001
$ws.Range("J113").Value = 1485.5
$ws.Range("K113").Value = 881.1429000000001
$ws.Range("L113").Value = 1485.5
$ws.Range("M113").Value = 1288.8571
$ws.Range("N113").Value = -5825.5
$ws.Range("H122").Value = 1171.2
$ws.Range("I122").Value = 1075
$ws.Range("J122").Value = 1235.3334
$ws.Range("K122").Value = 3225
$ws.Range("L122").Value = 3706.0002
$ws.Range("M122").Value = -775
$ws.Range("N122").Value = -8606.0002
$ws.Range("H132").Value = 3036.75
$ws.Range("I132").Value = 2924.7368
$ws.Range("J132").Value = 3462.4
$ws.Range("K132").Value = 8774.2104
$ws.Range("L132").Value = 10387.2
$ws.Range("M132").Value = -6244.2104
$ws.Range("N132").Value = -15447.2
$ws.Range("H136").Value = 1605.5897
$ws.Range("I136").Value = 1366.2667
$ws.Range("J136").Value = 2403.3333
$ws.Range("K136").Value = 4098.800099999999
$ws.Range("L136").Value = 7209.999899999999
$ws.Range("M136").Value = -1548.800099999999
$ws.Range("N136").Value = -12309.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 450
$ws.Range("J34").Value = 500
$ws.Range("L34").Value = 1500
$ws.Range("N34").Value = -1668
$ws.Range("H39").Value = 1833.3334
$ws.Range("H55").Value = 18796.666
$ws.Range("J55").Value = 3743.3333
$ws.Range("L55").Value = 11229.9999
$ws.Range("N55").Value = -11583.9999
$ws.Range("H113").Value = 1350.8667
$ws.Range("J113").Value = 769.2857
$ws.Range("L113").Value = 2307.8571
$ws.Range("N113").Value = -6647.8571
$ws.Range("H117").Value = 8343.75
$ws.Range("J117").Value = 8343.75
$ws.Range("L117").Value = 25031.25
$ws.Range("N117").Value = -31915.25
$ws.Range("H128").Value = 199000
$ws.Range("I128").Value = 199000
$ws.Range("K128").Value = 597000
$ws.Range("M128").Value = -592020

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 71589656
$ws.Range("I80").Value = 167034690
$ws.Range("J80").Value = 5893.75
$ws.Range("K80").Value = 167034690
$ws.Range("L80").Value = 5893.75
$ws.Range("M80").Value = -167033692
$ws.Range("N80").Value = -7889.75
$ws.Range("H83").Value = 71589656
$ws.Range("I83").Value = 167034690
$ws.Range("J83").Value = 5893.75
$ws.Range("K83").Value = 835173450
$ws.Range("L83").Value = 29468.75
$ws.Range("M83").Value = -835168458
$ws.Range("N83").Value = -39452.75
$ws.Range("H97").Value = 27778752
$ws.Range("I97").Value = 33334394
$ws.Range("J97").Value = 536
$ws.Range("K97").Value = 33334394
$ws.Range("L97").Value = 536
$ws.Range("M97").Value = -33333898
$ws.Range("N97").Value = -1528
$ws.Range("H132").Value = 7928.143
$ws.Range("I132").Value = 10750
$ws.Range("J132").Value = 6799.4
$ws.Range("K132").Value = 32250
$ws.Range("L132").Value = 20398.2
$ws.Range("M132").Value = -29720
$ws.Range("N132").Value = -25458.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4244.1333
$ws.Range("I61").Value = 3388.5
$ws.Range("K61").Value = 3388.5
$ws.Range("M61").Value = -3186.5
$ws.Range("H113").Value = 4244.1333
$ws.Range("I113").Value = 3388.5
$ws.Range("K113").Value = 3388.5
$ws.Range("M113").Value = -1218.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 62501332
$ws.Range("I96").Value = 111112696
$ws.Range("K96").Value = 111112696
$ws.Range("M96").Value = -111111323
$ws.Range("H132").Value = 30672.918
$ws.Range("I132").Value = 4195.727
$ws.Range("J132").Value = 69506.13
$ws.Range("K132").Value = 12587.181
$ws.Range("L132").Value = 208518.39
$ws.Range("M132").Value = -10057.181
$ws.Range("N132").Value = -213578.39
